$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Discard the previously-entered timestamp-looking values in column A.
$ws.Range("A3").ClearContents()
$ws.Range("A6").ClearContents()
$ws.Range("A9").ClearContents()

# Add the new "L:" / "R:" timing block starting at row 19.
$ws.Range("A19").Value = "L: "
$ws.Range("A20").Value = ":"
$ws.Range("A21").Value = 30

$ws.Range("A23").Value = ":"
$ws.Range("A24").Value = 31

$ws.Range("A26").Value = ":"
$ws.Range("A27").Value = 32

$ws.Range("A28").Value = "R:"
$ws.Range("A29").Value = ":"
$ws.Range("A30").Value = 10

$ws.Range("A32").Value = ":"
$ws.Range("A33").Value = 11

$ws.Range("A35").Value = ":"
$ws.Range("A36").Value = 12
